$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.470.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.093.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5221"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4435"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +15.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08945"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.087.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.691"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.712"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001122"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06610"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.272"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.507.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.316"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.331.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.184"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.651"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.157"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.901"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02557"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06794"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.465"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2255"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6890"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.251"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6326"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.201"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("E49").Value = "  +5.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.244"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.89%  "

Write-Output "Updated cryptos list"